$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Datos actualizados" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Marzo de 2020 a las 21:45"

# 2. Move the "Tenerife" row from its old position (row 30, right after
#    Pontevedra) up to right after "Murcia" (row 20, before Guadalajara),
#    and update its "Casos totales" (B) and "Muertes" (E) figures.
#    Insert a fresh row at 20 to make room; this pushes
#    Guadalajara..Pontevedra down by one (rows 20-29 -> 21-30), and the old
#    Tenerife row (still holding its stale data) shifts from row 30 to row 31.
$ws.Rows(20).Insert()
$ws.Range("A20").Value = "Tenerife"
$ws.Range("B20").Value = 210
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 135
$ws.Range("E20").Value = 3

# Remove the now-duplicate old Tenerife row (shifted down to row 31).
$ws.Rows(31).Delete()

# 3. Swap "Ceuta" and "La Palma" ordering (rows 57 and 58); their numeric
#    data (5, 0, 5, 0) is identical so only the labels need to swap.
$ws.Range("A57").Value = "Ceuta"
$ws.Range("A58").Value = "La Palma"
